$wb = $excel.ActiveWorkbook

# Rename sheets (sheet1..sheet5 map to GNG_TO, NB_TO, RS_TO, TOL_TO, vSAT_TO respectively)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16502912554544604"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16502912581383464"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16502912581383464"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16502912582010183"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1650291258280236"

# Sheet1 (GNG_TO) cell updates
$ws1.Range("B2").Value = "go_stims-16502912554246917.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912554379964.csv"
$ws1.Range("B4").Value = "go_stims-16502912554398026.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912554534283.csv"

# Sheet2 (NB_TO) cell updates
$ws2.Range("B2").Value = "ZB-match_0-16502912560713222.csv"
$ws2.Range("B3").Value = "ZB-match_2-16502912563351932.csv"
$ws2.Range("B4").Value = "TB-16502912581246355.csv"
$ws2.Range("B5").Value = "OB-16502912566635165.csv"
$ws2.Range("B6").Value = "TB-16502912578361993.csv"
$ws2.Range("B7").Value = "OB-1650291256517527.csv"
$ws2.Range("B8").Value = "OB-16502912571252935.csv"
$ws2.Range("B9").Value = "TB-16502912574955223.csv"
$ws2.Range("B10").Value = "ZB-match_6-16502912558641214.csv"

# Sheet4 (TOL_TO) cell updates
$ws4.Range("B2").Value = "MM_stims-16502912581542041.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912581415367.csv"
$ws4.Range("B4").Value = "MM_stims-16502912581853528.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912581552098.csv"
$ws4.Range("B6").Value = "MM_stims-16502912582010183.csv"
$ws4.Range("B7").Value = "ZM_stims-1650291258186356.csv"

# Sheet5 (vSAT_TO) cell updates
$ws5.Range("B2").Value = "SAT_stims-16502912582044687.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912582167177.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502912582623954.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912582322955.csv"
